# Applies the addition of rows 64-68 (new species-find records) to the
# "Artfynd" worksheet, matching the appended rows in the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 64
$ws.Cells.Item(64, 1).Value = 111998587
$ws.Cells.Item(64, 2).Value = 85192
$ws.Cells.Item(64, 3).Value = "Ovaliderad"
$ws.Cells.Item(64, 4).Value = "VU"
$ws.Cells.Item(64, 5).Value = 3595
$ws.Cells.Item(64, 6).Value = "Siljansspindling"
$ws.Cells.Item(64, 7).Value = "Cortinarius dalecarlicus"
$ws.Cells.Item(64, 8).Value = "Brandrud"
$ws.Cells.Item(64, 16).Value = "Bye, Jmt"
$ws.Cells.Item(64, 17).Value = 485438.9101868912
$ws.Cells.Item(64, 18).Value = 6995892.650789962
$ws.Cells.Item(64, 19).Value = 25
$ws.Cells.Item(64, 20).Value = "Jämtland"
$ws.Cells.Item(64, 21).Value = "Östersund"
$ws.Cells.Item(64, 22).Value = "Jämtland"
$ws.Cells.Item(64, 23).Value = "Marieby"
$ws.Cells.Item(64, 25).NumberFormat = "@"
$ws.Cells.Item(64, 25).Value = "2023-09-09"
$ws.Cells.Item(64, 25).Style = "Normal"
$ws.Cells.Item(64, 26).NumberFormat = "@"
$ws.Cells.Item(64, 26).Value = "16:38"
$ws.Cells.Item(64, 26).Style = "Normal"
$ws.Cells.Item(64, 27).NumberFormat = "@"
$ws.Cells.Item(64, 27).Value = "2023-09-09"
$ws.Cells.Item(64, 27).Style = "Normal"
$ws.Cells.Item(64, 28).NumberFormat = "@"
$ws.Cells.Item(64, 28).Value = "16:38"
$ws.Cells.Item(64, 28).Style = "Normal"
$ws.Cells.Item(64, 30).Value = $False
$ws.Cells.Item(64, 31).Value = $False
$ws.Cells.Item(64, 33).Value = $False
$ws.Cells.Item(64, 49).Value = "Johan Råghall"
$ws.Cells.Item(64, 50).Value = "Johan Råghall, Maria Danvind, Lars-Olof Grund, Magnus Andersson"

# Row 65
$ws.Cells.Item(65, 1).Value = 111998584
$ws.Cells.Item(65, 2).Value = 90655
$ws.Cells.Item(65, 3).Value = "Ovaliderad"
$ws.Cells.Item(65, 4).Value = "VU"
$ws.Cells.Item(65, 5).Value = 150
$ws.Cells.Item(65, 6).Value = "Grangråticka"
$ws.Cells.Item(65, 7).Value = "Boletopsis leucomelaena"
$ws.Cells.Item(65, 8).Value = "(Pers.) Fayod"
$ws.Cells.Item(65, 16).Value = "Bye, Jmt"
$ws.Cells.Item(65, 17).Value = 485432.9546544506
$ws.Cells.Item(65, 18).Value = 6995879.54152041
$ws.Cells.Item(65, 19).Value = 25
$ws.Cells.Item(65, 20).Value = "Jämtland"
$ws.Cells.Item(65, 21).Value = "Östersund"
$ws.Cells.Item(65, 22).Value = "Jämtland"
$ws.Cells.Item(65, 23).Value = "Marieby"
$ws.Cells.Item(65, 25).NumberFormat = "@"
$ws.Cells.Item(65, 25).Value = "2023-09-09"
$ws.Cells.Item(65, 25).Style = "Normal"
$ws.Cells.Item(65, 26).NumberFormat = "@"
$ws.Cells.Item(65, 26).Value = "16:53"
$ws.Cells.Item(65, 26).Style = "Normal"
$ws.Cells.Item(65, 27).NumberFormat = "@"
$ws.Cells.Item(65, 27).Value = "2023-09-09"
$ws.Cells.Item(65, 27).Style = "Normal"
$ws.Cells.Item(65, 28).NumberFormat = "@"
$ws.Cells.Item(65, 28).Value = "16:53"
$ws.Cells.Item(65, 28).Style = "Normal"
$ws.Cells.Item(65, 30).Value = $False
$ws.Cells.Item(65, 31).Value = $False
$ws.Cells.Item(65, 33).Value = $False
$ws.Cells.Item(65, 49).Value = "Johan Råghall"
$ws.Cells.Item(65, 50).Value = "Johan Råghall, Maria Danvind, Lars-Olof Grund, Magnus Andersson"

# Row 66
$ws.Cells.Item(66, 1).Value = 111998588
$ws.Cells.Item(66, 2).Value = 88955
$ws.Cells.Item(66, 3).Value = "Ovaliderad"
$ws.Cells.Item(66, 4).Value = "VU"
$ws.Cells.Item(66, 5).Value = 233196
$ws.Cells.Item(66, 6).Value = "Fjällfotad fingersvamp"
$ws.Cells.Item(66, 7).Value = "Ramaria rufescens"
$ws.Cells.Item(66, 8).Value = "(Schaeff.) Corner"
$ws.Cells.Item(66, 16).Value = "Bye, Jmt"
$ws.Cells.Item(66, 17).Value = 485478.8050299661
$ws.Cells.Item(66, 18).Value = 6995887.942324123
$ws.Cells.Item(66, 19).Value = 25
$ws.Cells.Item(66, 20).Value = "Jämtland"
$ws.Cells.Item(66, 21).Value = "Östersund"
$ws.Cells.Item(66, 22).Value = "Jämtland"
$ws.Cells.Item(66, 23).Value = "Marieby"
$ws.Cells.Item(66, 25).NumberFormat = "@"
$ws.Cells.Item(66, 25).Value = "2023-09-09"
$ws.Cells.Item(66, 25).Style = "Normal"
$ws.Cells.Item(66, 26).NumberFormat = "@"
$ws.Cells.Item(66, 26).Value = "16:31"
$ws.Cells.Item(66, 26).Style = "Normal"
$ws.Cells.Item(66, 27).NumberFormat = "@"
$ws.Cells.Item(66, 27).Value = "2023-09-09"
$ws.Cells.Item(66, 27).Style = "Normal"
$ws.Cells.Item(66, 28).NumberFormat = "@"
$ws.Cells.Item(66, 28).Value = "16:31"
$ws.Cells.Item(66, 28).Style = "Normal"
$ws.Cells.Item(66, 30).Value = $False
$ws.Cells.Item(66, 31).Value = $False
$ws.Cells.Item(66, 33).Value = $False
$ws.Cells.Item(66, 49).Value = "Johan Råghall"
$ws.Cells.Item(66, 50).Value = "Johan Råghall, Maria Danvind, Lars-Olof Grund, Magnus Andersson"

# Row 67
$ws.Cells.Item(67, 1).Value = 111998589
$ws.Cells.Item(67, 2).Value = 88950
$ws.Cells.Item(67, 3).Value = "Ovaliderad"
$ws.Cells.Item(67, 4).Value = "NT"
$ws.Cells.Item(67, 5).Value = 256756
$ws.Cells.Item(67, 6).Value = "Blek fingersvamp"
$ws.Cells.Item(67, 7).Value = "Ramaria pallida"
$ws.Cells.Item(67, 8).Value = "(Schaeff.) Ricken"
$ws.Cells.Item(67, 16).Value = "Bye, Jmt"
$ws.Cells.Item(67, 17).Value = 485478.7705635355
$ws.Cells.Item(67, 18).Value = 6995880.242057818
$ws.Cells.Item(67, 19).Value = 25
$ws.Cells.Item(67, 20).Value = "Jämtland"
$ws.Cells.Item(67, 21).Value = "Östersund"
$ws.Cells.Item(67, 22).Value = "Jämtland"
$ws.Cells.Item(67, 23).Value = "Marieby"
$ws.Cells.Item(67, 25).NumberFormat = "@"
$ws.Cells.Item(67, 25).Value = "2023-09-09"
$ws.Cells.Item(67, 25).Style = "Normal"
$ws.Cells.Item(67, 26).NumberFormat = "@"
$ws.Cells.Item(67, 26).Value = "16:20"
$ws.Cells.Item(67, 26).Style = "Normal"
$ws.Cells.Item(67, 27).NumberFormat = "@"
$ws.Cells.Item(67, 27).Value = "2023-09-09"
$ws.Cells.Item(67, 27).Style = "Normal"
$ws.Cells.Item(67, 28).NumberFormat = "@"
$ws.Cells.Item(67, 28).Value = "16:20"
$ws.Cells.Item(67, 28).Style = "Normal"
$ws.Cells.Item(67, 30).Value = $False
$ws.Cells.Item(67, 31).Value = $False
$ws.Cells.Item(67, 33).Value = $False
$ws.Cells.Item(67, 49).Value = "Johan Råghall"
$ws.Cells.Item(67, 50).Value = "Johan Råghall, Maria Danvind, Lars-Olof Grund, Magnus Andersson"

# Row 68
$ws.Cells.Item(68, 1).Value = 111998585
$ws.Cells.Item(68, 2).Value = 90669
$ws.Cells.Item(68, 3).Value = "Ovaliderad"
$ws.Cells.Item(68, 4).Value = "VU"
$ws.Cells.Item(68, 5).Value = 6003297
$ws.Cells.Item(68, 6).Value = "Spricktaggsvamp"
$ws.Cells.Item(68, 7).Value = "Hydnellum glaucopus"
$ws.Cells.Item(68, 8).Value = "(Maas Geest. & Nannf.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Cells.Item(68, 16).Value = "Bye, Jmt"
$ws.Cells.Item(68, 17).Value = 485448.7995749199
$ws.Cells.Item(68, 18).Value = 6995872.675907309
$ws.Cells.Item(68, 19).Value = 25
$ws.Cells.Item(68, 20).Value = "Jämtland"
$ws.Cells.Item(68, 21).Value = "Östersund"
$ws.Cells.Item(68, 22).Value = "Jämtland"
$ws.Cells.Item(68, 23).Value = "Marieby"
$ws.Cells.Item(68, 25).NumberFormat = "@"
$ws.Cells.Item(68, 25).Value = "2023-09-09"
$ws.Cells.Item(68, 25).Style = "Normal"
$ws.Cells.Item(68, 26).NumberFormat = "@"
$ws.Cells.Item(68, 26).Value = "16:48"
$ws.Cells.Item(68, 26).Style = "Normal"
$ws.Cells.Item(68, 27).NumberFormat = "@"
$ws.Cells.Item(68, 27).Value = "2023-09-09"
$ws.Cells.Item(68, 27).Style = "Normal"
$ws.Cells.Item(68, 28).NumberFormat = "@"
$ws.Cells.Item(68, 28).Value = "16:48"
$ws.Cells.Item(68, 28).Style = "Normal"
$ws.Cells.Item(68, 30).Value = $False
$ws.Cells.Item(68, 31).Value = $False
$ws.Cells.Item(68, 33).Value = $False
$ws.Cells.Item(68, 49).Value = "Johan Råghall"
$ws.Cells.Item(68, 50).Value = "Johan Råghall, Maria Danvind, Lars-Olof Grund, Magnus Andersson"

